$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.931.68'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.36%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.775.99'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.61%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.996'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.30%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.66'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.80'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.56%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.17%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.10%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.22%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.446'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.80%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.58'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +4.46%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000244'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.42%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '35.26'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.81%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.412.10'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.59%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.794.40'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.56%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.900.48'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.33%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.77%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.90%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.99'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.23%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '456.91'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.99%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.56'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.57%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.693'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '82.81'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.63%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -6.62%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.89'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.61%  '

$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("B27").Value = 'Fetch.AI'
$ws.Range("C27").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.07'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.40%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.87'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.28%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.929.00'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.27%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.26%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.21'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.47%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -7.61%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '28.91'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.26%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.41%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.90'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.84%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0987'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.40%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.99%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.77'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.14%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.90%  '

$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.09%  '

$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.14'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -6.81%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.00%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '47.24'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.78%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '43.14'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.54%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '152.25'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.94%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.294'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.27%  '

$ws.Range("B47").Value = 'Notcoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/2L2Y4ghjj+notcoin-not'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0271'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +90.09%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.04%  '

$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.28'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.60%  '

$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.84'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.13%  '

$ws.Range("B51").Value = 'Bittensor'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '384.91'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.52%  '
